# Agrupamento de distribuidoras: ENF + EMG = EMR, e EBO + EPB = EPB
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-RowByName($ws, $name, $lastRow) {
    for ($r = 2; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value2 -eq $name) {
            return $r
        }
    }
    return -1
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$rEBO = Find-RowByName $ws "EBO" $lastRow
$rEPB = Find-RowByName $ws "EPB" $lastRow
$rEMG = Find-RowByName $ws "EMG" $lastRow
$rENF = Find-RowByName $ws "ENF" $lastRow

# Sum EBO + EPB row-wise across columns B:I (new EPB totals)
$epbSums = @()
for ($c = 2; $c -le 9; $c++) {
    $epbSums += ($ws.Cells.Item($rEBO, $c).Value2 + $ws.Cells.Item($rEPB, $c).Value2)
}

# Sum EMG + ENF row-wise across columns B:I (new EMR totals)
$emrSums = @()
for ($c = 2; $c -le 9; $c++) {
    $emrSums += ($ws.Cells.Item($rEMG, $c).Value2 + $ws.Cells.Item($rENF, $c).Value2)
}

# Remove the four source rows (EBO, EMG, ENF fully disappear; EPB is removed here
# and re-created at the bottom of the table with its combined totals).
# Delete from the bottom-most row upward so earlier row numbers stay valid.
$rowsToDelete = @($rEBO, $rEPB, $rEMG, $rENF) | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Append the new EMR row, then the updated EPB row, as the last two rows of the table
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$emrRow = $lastRow + 1
$epbRow = $lastRow + 2

$ws.Cells.Item($emrRow, 1).Value = "EMR"
for ($c = 2; $c -le 9; $c++) {
    $ws.Cells.Item($emrRow, $c).Value = $emrSums[$c - 2]
}

$ws.Cells.Item($epbRow, 1).Value = "EPB"
for ($c = 2; $c -le 9; $c++) {
    $ws.Cells.Item($epbRow, $c).Value = $epbSums[$c - 2]
}

# Match the author's final on-screen selection
$ws.Range("L10").Select() | Out-Null
